# Insert 2 new data rows into the "Hortaliza, Vega Modelo de Temuco - Cebolla"
# sheet at row 1006 (pushing the existing rows 1006:1095 down to 1008:1097),
# then populate the two newly inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 1006 (shifts everything at/after row 1006
# down by two rows; formatting of row 1006 - the date style on column D - is
# carried down from the row above, matching the source workbook's style).
$ws.Range("A1006:A1007").EntireRow.Insert()

# New row 1006
$ws.Range("A1006").Value = 10
$ws.Range("B1006").Value = "Vega Modelo de Temuco"
$ws.Range("C1006").Value = "La Araucanía"
$ws.Range("D1006").Value = 44769
$ws.Range("E1006").Value = 9
$ws.Range("F1006").Value = 100112004
$ws.Range("G1006").Value = "Cebolla"
$ws.Range("H1006").Value = "Morada(o)"
$ws.Range("I1006").Value = "Primera"
$ws.Range("J1006").Value = 100
$ws.Range("K1006").Value = 16000
$ws.Range("L1006").Value = 16000
$ws.Range("M1006").Value = 16000
$ws.Range("N1006").Value = "$/malla 18 kilos"
$ws.Range("O1006").Value = "Perú"
$ws.Range("P1006").Value = 889
$ws.Range("Q1006").Value = 18
$ws.Range("R1006").Value = "Hortaliza"

# New row 1007
$ws.Range("A1007").Value = 10
$ws.Range("B1007").Value = "Vega Modelo de Temuco"
$ws.Range("C1007").Value = "La Araucanía"
$ws.Range("D1007").Value = 44769
$ws.Range("E1007").Value = 9
$ws.Range("F1007").Value = 100112004
$ws.Range("G1007").Value = "Cebolla"
$ws.Range("H1007").Value = "Sin especificar"
$ws.Range("I1007").Value = "1a (guarda)"
$ws.Range("J1007").Value = 500
$ws.Range("K1007").Value = 6000
$ws.Range("L1007").Value = 7000
$ws.Range("M1007").Value = 6600
$ws.Range("N1007").Value = "$/malla 18 kilos"
$ws.Range("O1007").Value = "Región de O'Higgins"
$ws.Range("P1007").Value = 367
$ws.Range("Q1007").Value = 18
$ws.Range("R1007").Value = "Hortaliza"
